# "remove duplicated plugins and update browser capabilities"
#
# Sheet "c-demo_ipad" (sheet1) has a key/value config table. The webdriver
# browser-capability rows (6-12) contained duplicated "platform"/"version"
# keys for chrome/ie/firefox. They are renamed/replaced with the current
# msedgedriver + platformName/browserVersion style keys, and the selected
# cell in the sheet view moves to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("c-demo_ipad")

# Row 6: webdriver.path.ie -> webdriver.path.edge
$ws.Range("A6").Value = "webdriver.path.edge"
$ws.Range("B6").Value = "web_drivers/windows/msedgedriver.exe"

# Row 7: webdriver.platform.chrome -> webdriver.platformName.chrome
$ws.Range("A7").Value = "webdriver.platformName.chrome"

# Row 8: webdriver.version.chrome -> webdriver.browserVersion.chrome (61.0 -> 91.0)
$ws.Range("A8").Value = "webdriver.browserVersion.chrome"
$ws.Range("B8").Value = "91.0"

# Row 9: webdriver.platform.ie -> webdriver.platformName.edge
$ws.Range("A9").Value = "webdriver.platformName.edge"

# Row 10: webdriver.version.ie -> webdriver.browserVersion.edge (11 -> 91)
$ws.Range("A10").Value = "webdriver.browserVersion.edge"
$ws.Range("B10").Value = "91"

# Row 11: webdriver.platform.firefox -> webdriver.platformName.firefox
$ws.Range("A11").Value = "webdriver.platformName.firefox"

# Row 12: webdriver.version.firefox -> webdriver.browserVersion.firefox (60 -> 96)
$ws.Range("A12").Value = "webdriver.browserVersion.firefox"
$ws.Range("B12").Value = "96"

# Update the active selection to reflect where the editor left off.
$ws.Range("A7").Select()
